$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix F12's alignment to match the other time cells in the column (center).
$ws.Range("F12").HorizontalAlignment = -4108

# Add new row 13: "Lambda Functions" task covering Edge customization topics.
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = $ws.Range("D12").Value2
$ws.Range("E13").Value = "customization at Edge, CloudFront Func, Lambda@Edge, Lambda by default, VPC, lambda layers."
$ws.Range("F13").Value = 0.70833333333333337
$ws.Range("G13").Value = 0.75416666666666665

# Match formatting of the row above for the new row (same styles as row 12,
# i.e. G12's un-centered time style for both F13 and G13).
$ws.Range("D13").HorizontalAlignment = $ws.Range("D12").HorizontalAlignment
$ws.Range("C13").HorizontalAlignment = $ws.Range("C12").HorizontalAlignment
$ws.Range("F13").NumberFormat = $ws.Range("G12").NumberFormat
$ws.Range("G13").NumberFormat = $ws.Range("G12").NumberFormat
